# Clean whitespace from country column and update values per source data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country names (trimmed of leading/trailing whitespace) and updated numeric data
$data = @(
    @{ Row = 2;  Country = "United Kingdom"; B = 60; C = 237 },
    @{ Row = 3;  Country = "Greece";         B = 18; C = 734 },
    @{ Row = 4;  Country = "United States";  B = 18; C = 158 },
    @{ Row = 5;  Country = "France";         B = 16; C = 397 },
    @{ Row = 6;  Country = "Italy";          B = 15; C = 713 },
    @{ Row = 7;  Country = "Norway";         B = 8;  C = 178 },
    @{ Row = 8;  Country = "Poland";         B = 8;  C = 90 },
    @{ Row = 9;  Country = "Iceland";        B = 7;  C = 25 },
    @{ Row = 10; Country = "Japan";          B = 7;  C = 303 },
    @{ Row = 11; Country = "Romania";        B = 7;  C = 70 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Country
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
}
